$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Top "Bad Drivers" table -------------------------------------------------
$ws.Range("C3").Value = 7
$ws.Range("D3").Value = 98.5
$ws.Range("C4").Value = 7

# --- "Good Drivers" table (rows 12-17) reordered / updated ------------------
# Row 12
$ws.Range("A12").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B12").Value = 445055
$ws.Range("D12").Value = 99.90000000000001
# Leading apostrophe keeps this as literal text "2024-11-10" instead of
# having Excel auto-convert the date-shaped string into a date serial.
$ws.Range("E12").Value = "'2024-11-10"

# Row 13
$ws.Range("A13").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B13").Value = 77849
$ws.Range("D13").Value = 99.90000000000001
$ws.Range("E13").Value = "'2021-08-18"

# Row 14
$ws.Range("A14").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B14").Value = 34244
$ws.Range("D14").Value = 100
$ws.Range("E14").Value = "'2021-04-27"

# Row 15
$ws.Range("A15").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B15").Value = 59673
$ws.Range("D15").Value = 100
$ws.Range("E15").Value = "'2020-08-05"

# Row 16 - D16 (100) is unchanged by the diff, so it is left alone.
$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B16").Value = 113652
$ws.Range("E16").Value = "'2020-01-06"

# Row 17 - only the adapter name and client count move here; D17/E17 (100,
# "2019-12-14") are already correct and untouched by the diff.
$ws.Range("A17").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B17").Value = 56018
